# Apply updated symbol/price list (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.06'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.89'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.272'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05731'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.441'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8103'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8791'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1445'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07369'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03124'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09394'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001598'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04804'
$ws.Range("E15").Value = '14CoinExTokenCET'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0005850'
$ws.Range("E16").Value = '15OneONE'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006152'
$ws.Range("E17").Value = '16TigerCashTCH'
$ws.Range("B18").Value = 'HotbitToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005100'
$ws.Range("E18").Value = '17HotbitTokenHTB'
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0009974'
$ws.Range("E19").Value = '18BitKanKAN'
$ws.Range("B20").Value = 'NitroEx'
$ws.Range("C20").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0001500'
$ws.Range("E20").Value = '19NitroExNTX'
$ws.Range("B21").Value = 'LEO'
$ws.Range("C21").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.749'
$ws.Range("E21").Value = '20LEOLEO'
$ws.Range("B22").Value = 'KuCoinToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.302'
$ws.Range("E22").Value = '21KuCoinTokenKCS'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.186'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'BitpandaEcosystemToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.3279'
$ws.Range("E24").Value = '23BitpandaEcosystemTokenBEST'
$ws.Range("B25").Value = 'ProBitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1320'
$ws.Range("E25").Value = '24ProBitTokenPROB'
$ws.Range("B26").Value = 'MCDex'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.134'
$ws.Range("E26").Value = '25MCDexMCB'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003001'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03906'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006770'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003200'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007477'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005638'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6000'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1739'
